$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93
$ws.Range("B93").Value = 6236611
$ws.Range("E93").Value = 'Mineros'
$ws.Range("F93").Value = 'Monagas'
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 1
$ws.Range("I93").Value = 'H'
$ws.Range("J93").Value = 3.2
$ws.Range("K93").Value = 3.4
$ws.Range("L93").Value = 2
$ws.Range("M93").Value = 4.2
$ws.Range("N93").Value = 3.8
$ws.Range("O93").Value = 1.65
$ws.Range("P93").Value = 0.75
$ws.Range("Q93").Value = 1.95
$ws.Range("R93").Value = 1.85
$ws.Range("S93").Value = 2.5
$ws.Range("T93").Value = 1.825
$ws.Range("U93").Value = 1.975
$ws.Range("V93").Value = 3.2
$ws.Range("W93").Value = -1
$ws.Range("X93").Value = -1
$ws.Range("Y93").Value = 0.95
$ws.Range("Z93").Value = -1
$ws.Range("AA93").Value = 0.825
$ws.Range("AB93").Value = -1

# Row 94
$ws.Range("B94").Value = 6236255
$ws.Range("E94").Value = 'Deportivo Rayo Zuliano'
$ws.Range("F94").Value = 'Caracas'
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 'D'
$ws.Range("J94").Value = 3.75
$ws.Range("K94").Value = 3.1
$ws.Range("L94").Value = 1.95
$ws.Range("M94").Value = 2.9
$ws.Range("N94").Value = 2.875
$ws.Range("O94").Value = 2.45
$ws.Range("P94").Value = 0.25
$ws.Range("Q94").Value = 1.775
$ws.Range("R94").Value = 2.025
$ws.Range("S94").Value = 2.25
$ws.Range("T94").Value = 1.85
$ws.Range("U94").Value = 1.95
$ws.Range("V94").Value = -1
$ws.Range("W94").Value = 1.875
$ws.Range("X94").Value = -1
$ws.Range("Y94").Value = 0.3875
$ws.Range("Z94").Value = -0.5
$ws.Range("AA94").Value = -1
$ws.Range("AB94").Value = 0.95

# Row 95
$ws.Range("B95").Value = 6236251
$ws.Range("E95").Value = 'Angostura FC'
$ws.Range("F95").Value = 'Portuguesa'
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 2
$ws.Range("I95").Value = 'A'
$ws.Range("J95").Value = 3.1
$ws.Range("K95").Value = 3.2
$ws.Range("L95").Value = 2.15
$ws.Range("M95").Value = 4
$ws.Range("N95").Value = 3.6
$ws.Range("O95").Value = 1.75
$ws.Range("P95").Value = 0.75
$ws.Range("Q95").Value = 1.8
$ws.Range("R95").Value = 2
$ws.Range("S95").Value = 2.5
$ws.Range("T95").Value = 1.95
$ws.Range("U95").Value = 1.85
$ws.Range("V95").Value = -1
$ws.Range("W95").Value = -1
$ws.Range("X95").Value = 0.75
$ws.Range("Y95").Value = -0.5
$ws.Range("Z95").Value = 0.5
$ws.Range("AA95").Value = 0.95
$ws.Range("AB95").Value = -1

# Row 96
$ws.Range("B96").Value = 6236254
$ws.Range("E96").Value = 'Academia Puerto Cabello'
$ws.Range("F96").Value = 'Estudiantes Merida'
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 'H'
$ws.Range("J96").Value = 1.727
$ws.Range("K96").Value = 3.4
$ws.Range("L96").Value = 4.333
$ws.Range("M96").Value = 1.666
$ws.Range("N96").Value = 3.4
$ws.Range("O96").Value = 4.75
$ws.Range("P96").Value = -0.75
$ws.Range("Q96").Value = 1.875
$ws.Range("R96").Value = 1.925
$ws.Range("S96").Value = 2.5
$ws.Range("T96").Value = 1.9
$ws.Range("U96").Value = 1.9
$ws.Range("V96").Value = 0.6659999999999999
$ws.Range("W96").Value = -1
$ws.Range("X96").Value = -1
$ws.Range("Y96").Value = 0.4375
$ws.Range("Z96").Value = -0.5
$ws.Range("AA96").Value = -1
$ws.Range("AB96").Value = 0.8999999999999999

# Row 97
$ws.Range("B97").Value = 6236612
$ws.Range("E97").Value = 'Zamora'
$ws.Range("F97").Value = 'Carabobo'
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 2
$ws.Range("I97").Value = 'A'
$ws.Range("J97").Value = 3.2
$ws.Range("K97").Value = 3.1
$ws.Range("L97").Value = 2.15
$ws.Range("M97").Value = 4.5
$ws.Range("N97").Value = 3.3
$ws.Range("O97").Value = 1.75
$ws.Range("P97").Value = 0.5
$ws.Range("Q97").Value = 2
$ws.Range("R97").Value = 1.8
$ws.Range("S97").Value = 2.25
$ws.Range("T97").Value = 1.925
$ws.Range("U97").Value = 1.875
$ws.Range("V97").Value = -1
$ws.Range("W97").Value = -1
$ws.Range("X97").Value = 0.75
$ws.Range("Y97").Value = -1
$ws.Range("Z97").Value = 0.8
$ws.Range("AA97").Value = -0.5
$ws.Range("AB97").Value = 0.4375

# Row 157
$ws.Range("B157").Value = 7920997
$ws.Range("E157").Value = 'Carabobo'
$ws.Range("F157").Value = 'UCV'
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 1
$ws.Range("I157").Value = 'A'
$ws.Range("J157").Value = 1.833
$ws.Range("K157").Value = 3.1
$ws.Range("L157").Value = 4.2
$ws.Range("M157").Value = 1.833
$ws.Range("N157").Value = 3.1
$ws.Range("O157").Value = 4.2
$ws.Range("P157").Value = -0.5
$ws.Range("Q157").Value = 1.9
$ws.Range("R157").Value = 1.9
$ws.Range("S157").Value = 2
$ws.Range("T157").Value = 1.85
$ws.Range("U157").Value = 1.95
$ws.Range("V157").Value = -1
$ws.Range("W157").Value = -1
$ws.Range("X157").Value = 3.2
$ws.Range("Y157").Value = -1
$ws.Range("Z157").Value = 0.8999999999999999
$ws.Range("AA157").Value = -1
$ws.Range("AB157").Value = 0.95

# Row 158
$ws.Range("B158").Value = 7920998
$ws.Range("E158").Value = 'Zamora'
$ws.Range("F158").Value = 'Caracas'
$ws.Range("G158").Value = 2
$ws.Range("H158").Value = 2
$ws.Range("I158").Value = 'D'
$ws.Range("J158").Value = 3.75
$ws.Range("K158").Value = 3.2
$ws.Range("L158").Value = 1.909
$ws.Range("M158").Value = 3
$ws.Range("N158").Value = 2.9
$ws.Range("O158").Value = 2.375
$ws.Range("P158").Value = 0.25
$ws.Range("Q158").Value = 1.8
$ws.Range("R158").Value = 2
$ws.Range("S158").Value = 2
$ws.Range("T158").Value = 1.825
$ws.Range("U158").Value = 1.975
$ws.Range("V158").Value = -1
$ws.Range("W158").Value = 1.9
$ws.Range("X158").Value = -1
$ws.Range("Y158").Value = 0.4
$ws.Range("Z158").Value = -0.5
$ws.Range("AA158").Value = 0.825
$ws.Range("AB158").Value = -1
